$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.993.09"
$ws.Range("E2").Value = "  -6.71%  "
$ws.Range("D3").Value = "2.480.10"
$ws.Range("E3").Value = "  -12.14%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "471.67"
$ws.Range("E5").Value = "  -6.23%  "
$ws.Range("D6").Value = "134.35"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "0.495"
$ws.Range("E8").Value = "  -6.51%  "
$ws.Range("D9").Value = "2.499.64"
$ws.Range("E9").Value = "  -11.28%  "
$ws.Range("D10").Value = "0.0981"
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("D11").Value = "5.47"
$ws.Range("E11").Value = "  -7.56%  "
$ws.Range("D12").Value = "0.323"
$ws.Range("E12").Value = "  -6.93%  "
$ws.Range("E13").Value = "  -3.54%  "
$ws.Range("D14").Value = "2.915.42"
$ws.Range("E14").Value = "  -11.98%  "
$ws.Range("D15").Value = "54.891.29"
$ws.Range("E15").Value = "  -7.16%  "
$ws.Range("D16").Value = "0.0000138"
$ws.Range("E16").Value = "  +3.03%  "
$ws.Range("D17").Value = "20.15"
$ws.Range("E17").Value = "  -6.12%  "
$ws.Range("D18").Value = "2.481.75"
$ws.Range("E18").Value = "  -12.23%  "
$ws.Range("D19").Value = "4.26"
$ws.Range("E19").Value = "  -9.39%  "
$ws.Range("D20").Value = "316.52"
$ws.Range("E20").Value = "  -8.76%  "
$ws.Range("D21").Value = "9.65"
$ws.Range("E21").Value = "  -12.01%  "
$ws.Range("D22").Value = "1.01"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("D23").Value = "5.70"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").Value = "5.45"
$ws.Range("E24").Value = "  -12.76%  "
$ws.Range("D25").Value = "57.52"
$ws.Range("E25").Value = "  -8.58%  "
$ws.Range("E26").Value = "  +1.40%  "
$ws.Range("D27").Value = "0.390"
$ws.Range("E27").Value = "  -8.17%  "
$ws.Range("D28").Value = "0.158"
$ws.Range("E28").Value = "  -8.03%  "
$ws.Range("D29").Value = "2.536.04"
$ws.Range("E29").Value = "  -14.01%  "
$ws.Range("D30").Value = "7.34"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("D32").Value = "0.0₃0755"
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("D33").Value = "153.14"
$ws.Range("E33").Value = "  +2.01%  "
$ws.Range("D34").Value = "17.97"
$ws.Range("E34").Value = "  -4.88%  "
$ws.Range("D35").Value = "1.44"
$ws.Range("E35").Value = "  -9.94%  "
$ws.Range("D36").Value = "5.09"
$ws.Range("E36").Value = "  -3.87%  "
$ws.Range("D37").Value = "3.65"
$ws.Range("E37").Value = "  -11.55%  "
$ws.Range("D38").Value = "1.08"
$ws.Range("E38").Value = "  -3.38%  "
$ws.Range("D39").Value = "0.825"
$ws.Range("E39").Value = "  -7.54%  "
$ws.Range("D40").Value = "33.87"
$ws.Range("E40").Value = "  -7.54%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").Value = "0.614"
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "0.993"
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").Value = "0.0540"
$ws.Range("E43").Value = "  -2.68%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "3.34"
$ws.Range("E44").Value = "  -4.69%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "1.27"
$ws.Range("E45").Value = "  -5.30%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "10.21"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").Value = "1.978.06"
$ws.Range("E47").Value = "  -10.67%  "
$ws.Range("D48").Value = "0.0224"
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "4.54"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.0879"
$ws.Range("E50").Value = "  -0.88%  "
$ws.Range("D51").Value = "16.94"
$ws.Range("E51").Value = "  -11.26%  "
